$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 335, shifting the existing rows 335-354 down to 337-356.
$ws.Rows("335:336").Insert(-4121)

# New row 335: Choclero, "Choclero" priced entry for 2022-01-24 (serial 44585)
$ws.Range("A335").Value = 10
$ws.Range("B335").Value = "Vega Modelo de Temuco"
$ws.Range("C335").Value = "La Araucanía"
$ws.Range("D335").Value = 44585
$ws.Range("E335").Value = 9
$ws.Range("F335").Value = 100112024
$ws.Range("G335").Value = "Choclo"
$ws.Range("H335").Value = "Choclero"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 16000
$ws.Range("K335").Value = 200
$ws.Range("L335").Value = 250
$ws.Range("M335").Value = 222
$ws.Range("N335").Value = "$/unidad"
$ws.Range("O335").Value = "Región del Maule"
$ws.Range("P335").Value = 222
$ws.Range("Q335").Value = 1
$ws.Range("R335").Value = "Hortaliza"

# New row 336: Dulce o Americano priced entry for 2022-01-24 (serial 44585)
$ws.Range("A336").Value = 10
$ws.Range("B336").Value = "Vega Modelo de Temuco"
$ws.Range("C336").Value = "La Araucanía"
$ws.Range("D336").Value = 44585
$ws.Range("E336").Value = 9
$ws.Range("F336").Value = 100112024
$ws.Range("G336").Value = "Choclo"
$ws.Range("H336").Value = "Dulce o Americano"
$ws.Range("I336").Value = "Primera"
$ws.Range("J336").Value = 30000
$ws.Range("K336").Value = 200
$ws.Range("L336").Value = 200
$ws.Range("M336").Value = 200
$ws.Range("N336").Value = "$/unidad"
$ws.Range("O336").Value = "Región del Maule"
$ws.Range("P336").Value = 200
$ws.Range("Q336").Value = 1
$ws.Range("R336").Value = "Hortaliza"
